$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update C5 value (affects E5, F5, G5, M5, B6, C6, E6 via formulas)
$ws.Range("C5").Value = 1.0840000000000001

# Update B11 value (affects B12, B13, M6 via formulas)
$ws.Range("B11").Value = 3.1749999999999998

# New rows for plate and clamp
$ws.Range("A20").Value = "plate"
$ws.Range("B20").Value = 19
$ws.Range("J20").Value = 257
$ws.Range("L20").Formula = "=24*48 - 257*2 - 42.5*60*8"

$ws.Range("A21").Value = "clamp"
$ws.Range("B21").Value = 55.62
$ws.Range("J21").Value = 252

$ws.Range("B22").Formula = "=SUM(B20:B21)"

$ws.Range("L23").Formula = "=609.6 * 609.6*2 - 165806*2 - 42.5*60*8"

# M8 now references B22 (sum of plate+clamp) instead of B19
$ws.Range("M8").Formula = "=B22*8"

# Update selection to match target
$ws.Range("L23").Select()
